$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2-5
$ws.Range("B2").Value = 123

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 87

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 79

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 9

# Delete row 6 entirely (shifts cells up, removing the last row)
$ws.Range("A6:B6").Delete()
